# Regenerate s_val data to filter save games.
# Updates the numeric values in columns B:G for rows 2-6 of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.6545652718822623, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 3.536033448013082)
    3 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 5.582307763322248)
    4 = @(0.6545652718822623, 9.983522426115931, 0.7210945179870265, 13.86384647080068, 0, 25.2230286867859)
    5 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 13.86384647080068, 1, 18.91276827552123)
    6 = @(0.04172184405617529, 0.04103571897497393, 0.1496068669990043, 0.5333859586016987, 1, 0.7657503886318522)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2  # Column B is index 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
